$wb = $excel.ActiveWorkbook

# Sheet 1: "展览" (Exhibition)
$wsExhibition = $wb.Worksheets.Item(1)
$wsExhibition.Cells.Item(2, 6).Value = 35   # was 34
$wsExhibition.Cells.Item(3, 6).Value = 21172   # was 21156
$wsExhibition.Cells.Item(5, 6).Value = 339   # was 338
$wsExhibition.Cells.Item(8, 6).Value = 7898   # was 7891
$wsExhibition.Cells.Item(10, 6).Value = 41   # was 40
$wsExhibition.Cells.Item(11, 6).Value = 763   # was 762
$wsExhibition.Cells.Item(12, 6).Value = 309   # was 308
$wsExhibition.Cells.Item(19, 6).Value = 1359   # was 1357
$wsExhibition.Cells.Item(20, 6).Value = 523   # was 519
$wsExhibition.Cells.Item(25, 6).Value = 82   # was 81
$wsExhibition.Cells.Item(27, 6).Value = 1177   # was 1176
$wsExhibition.Cells.Item(33, 6).Value = 5   # was 4
$wsExhibition.Cells.Item(35, 6).Value = 5029   # was 5025
$wsExhibition.Cells.Item(36, 6).Value = 32   # was 31
$wsExhibition.Cells.Item(37, 6).Value = 99   # was 98
$wsExhibition.Cells.Item(40, 6).Value = 13056   # was 13050
$wsExhibition.Cells.Item(42, 6).Value = 129   # was 128
$wsExhibition.Cells.Item(46, 6).Value = 425   # was 423

# Sheet 4: "全部类型" (All Types)
$wsAllTypes = $wb.Worksheets.Item(4)
$wsAllTypes.Cells.Item(2, 6).Value = 35   # was 34
$wsAllTypes.Cells.Item(3, 6).Value = 21172   # was 21156
$wsAllTypes.Cells.Item(7, 6).Value = 7898   # was 7891
$wsAllTypes.Cells.Item(9, 6).Value = 41   # was 40
$wsAllTypes.Cells.Item(10, 6).Value = 763   # was 762
$wsAllTypes.Cells.Item(11, 6).Value = 309   # was 308
$wsAllTypes.Cells.Item(17, 6).Value = 1359   # was 1357
$wsAllTypes.Cells.Item(18, 6).Value = 523   # was 519
$wsAllTypes.Cells.Item(23, 6).Value = 82   # was 81
$wsAllTypes.Cells.Item(25, 6).Value = 1177   # was 1176
$wsAllTypes.Cells.Item(32, 6).Value = 5   # was 4
$wsAllTypes.Cells.Item(35, 6).Value = 5029   # was 5025
$wsAllTypes.Cells.Item(36, 6).Value = 32   # was 31
$wsAllTypes.Cells.Item(37, 6).Value = 99   # was 98
$wsAllTypes.Cells.Item(40, 6).Value = 13056   # was 13050
$wsAllTypes.Cells.Item(42, 6).Value = 129   # was 128
$wsAllTypes.Cells.Item(46, 6).Value = 425   # was 423
